# [UPDATE] Budget Return, All, Years and Duplicate
# Set all rows (2-49) in column A of the "rkap" sheet to year 2021
# (previously a mix of 2023 / 2022), and reset the sheet view (scroll
# position + selection) back to the top-left / B4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rkap")

# Update the "Year" column (A) for every data row (2 through 49) to 2021.
for ($r = 2; $r -le 49; $r++) {
    $ws.Cells.Item($r, 1).Value = 2021
}

# Activate the sheet, scroll back to the top-left corner and move the
# selection to B4 (previously topLeftCell="C1" / selection F13).
[void]$ws.Activate()
[void]$ws.Range("A1").Select()
[void]$ws.Range("B4").Select()
